$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.647240480131028
$ws.Range("K2").Value = 0.665933688883394
$ws.Range("L2").Value = 0.572690723237937
$ws.Range("N2").Value = 0.603471262597791

$ws.Range("B3").Value = 0.672054760884641
$ws.Range("C3").Value = 0.73150166637206
$ws.Range("K3").Value = 0.67292535125419
$ws.Range("L3").Value = 0.639033952437169
$ws.Range("N3").Value = 0.589468339412676

$ws.Range("B4").Value = 0.68518253119246
$ws.Range("K4").Value = 0.708149579691965
$ws.Range("L4").Value = 0.604866837966034
$ws.Range("N4").Value = 0.636828257869627

$ws.Range("B5").Value = 0.61895417279981
$ws.Range("K5").Value = 0.633861854984241
$ws.Range("L5").Value = 0.5411746546684
$ws.Range("N5").Value = 0.560767240057917
